$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (string) data type instead of
# being auto-converted to numbers/percentages by setting the cell format
# to Text ("@") before assigning the new value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

$ws.Range("D2").Value = "293.27"
$ws.Range("E2").Value = "0.39%"
$ws.Range("D3").Value = "40.49"
$ws.Range("E3").Value = "1.42%"
$ws.Range("D4").Value = "5.008"
$ws.Range("E4").Value = "-0.51%"
$ws.Range("D5").Value = "0.07358"
$ws.Range("E5").Value = "-0.27%"
$ws.Range("D6").Value = "1.568"
$ws.Range("E6").Value = "3.39%"
$ws.Range("E7").Value = "0.33%"
$ws.Range("D8").Value = "2.352"
$ws.Range("E8").Value = "-1.96%"
$ws.Range("D9").Value = "0.1167"
$ws.Range("E9").Value = "0.68%"
$ws.Range("D10").Value = "0.1778"
$ws.Range("E10").Value = "1.73%"
$ws.Range("D11").Value = "0.04382"
$ws.Range("E11").Value = "4.91%"
$ws.Range("D12").Value = "0.08753"
$ws.Range("E12").Value = "0.80%"
$ws.Range("D13").Value = "0.1055"
$ws.Range("E13").Value = "0.28%"
$ws.Range("D14").Value = "0.001266"
$ws.Range("E14").Value = "1.25%"
$ws.Range("D15").Value = "0.005916"
$ws.Range("E15").Value = "-1.53%"
$ws.Range("D16").Value = "3.350"
$ws.Range("E16").Value = "-0.18%"
$ws.Range("D17").Value = "4.292"
$ws.Range("E17").Value = "-0.39%"
$ws.Range("D19").Value = "7.822"
$ws.Range("E19").Value = "3.12%"
$ws.Range("D20").Value = "0.1389"
$ws.Range("E20").Value = "2.31%"
$ws.Range("E21").Value = "-1.69%"
$ws.Range("D22").Value = "0.03916"
$ws.Range("E22").Value = "2.06%"
$ws.Range("E23").Value = "-2.12%"
$ws.Range("D24").Value = "0.003676"
$ws.Range("E24").Value = "1.93%"
$ws.Range("E25").Value = "-8.30%"
$ws.Range("E26").Value = "-0.60%"
$ws.Range("D38").Value = "0.02342"
$ws.Range("E38").Value = "1.09%"
$ws.Range("D39").Value = "0.05098"
$ws.Range("E39").Value = "2.04%"
$ws.Range("D41").Value = "0.007854"
$ws.Range("E41").Value = "1.56%"
$ws.Range("E42").Value = "1.53%"
$ws.Range("D43").Value = "0.007383"
$ws.Range("E43").Value = "-0.69%"
$ws.Range("D44").Value = "0.008058"
$ws.Range("E44").Value = "1.69%"
$ws.Range("D45").Value = "0.2914"
$ws.Range("E45").Value = "-8.20%"
$ws.Range("D46").Value = "0.00006229"
$ws.Range("E46").Value = "-4.02%"
$ws.Range("E47").Value = "-0.59%"
$ws.Range("D48").Value = "0.04839"
$ws.Range("E48").Value = "-80.78%"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "-0.59%"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "-0.59%"
